# Update cryptocurrency price/volume data for Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.840.82"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").Value = "'1.751.93"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'333.52"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.3872"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("D8").Value = "'0.3383"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("D9").Value = "'45.43"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("D11").Value = "'0.07207"
$ws.Range("E11").Value = "  -2.34%  "

$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'22.51"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").Value = "'6.186"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.757.94"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.099"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "'0.00001059"
$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("D18").Value = "'0.06598"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").Value = "'79.49"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  -2.64%  "

$ws.Range("D22").Value = "'6.183"
$ws.Range("E22").Value = "  -3.08%  "

$ws.Range("D23").Value = "'27.869.90"
$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("E24").Value = "  -2.56%  "

$ws.Range("D25").Value = "'2.396"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").Value = "'154.43"
$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D27").Value = "'19.92"
$ws.Range("E27").Value = "  -3.06%  "

$ws.Range("D28").Value = "'2.304"
$ws.Range("E28").Value = "  -3.86%  "

$ws.Range("D29").Value = "'1.955.67"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").Value = "'1.285"
$ws.Range("E30").Value = "  -9.45%  "

$ws.Range("D31").Value = "'130.83"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "'4.021"
$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("D33").Value = "'5.822"
$ws.Range("E33").Value = "  -4.02%  "

$ws.Range("D34").Value = "'0.08791"
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("D35").Value = "'12.17"
$ws.Range("E35").Value = "  -3.71%  "

$ws.Range("D36").Value = "'1.538"
$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("D37").Value = "'0.6545"
$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("D38").Value = "'5.143"
$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").Value = "'0.02273"
$ws.Range("E39").Value = "  -5.44%  "

$ws.Range("D40").Value = "'0.06112"
$ws.Range("E40").Value = "  -2.50%  "

$ws.Range("D41").Value = "'0.2106"
$ws.Range("E41").Value = "  -2.96%  "

$ws.Range("D42").Value = "'1.210"
$ws.Range("E42").Value = "  -3.13%  "

$ws.Range("D43").Value = "'8.021"
$ws.Range("E43").Value = "  -1.96%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'13.64"
$ws.Range("E45").Value = "  -2.89%  "

$ws.Range("D46").Value = "'3.817"
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("D47").Value = "'0.6045"
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").Value = "'127.09"
$ws.Range("E48").Value = "  -2.88%  "

$ws.Range("D49").Value = "'1.994"
$ws.Range("E49").Value = "  -3.11%  "

$ws.Range("D50").Value = "'1.111"
$ws.Range("E50").Value = "  +4.87%  "

$ws.Range("D51").Value = "'1.163"
$ws.Range("E51").Value = "  +2.09%  "
